# NIT-9007250813.xlsx - "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" list (column E, rows 16-31) is re-sorted from descending
# (2104 ... 2001) to ascending (2001 ... 2104) order, and the "Valor Mora"
# (column F) travels together with its period: period 2104 keeps its 30506
# value (now on the last row, 31) while the remaining periods keep 35200
# (now with period 2001 on the first row, 16).
#
# We update the text/value content directly (rather than doing a generic
# Range.Sort, which would also drag the per-row cell formatting/styles along
# with it) so the existing borders/number-formats stay anchored to their
# original rows, exactly like the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2001", "2002", "2003", "2004", "2005", "2006", "2007", "2008", "2009", "2010", "2011", "2012", "2101", "2102", "2103", "2104")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# Valor Mora follows the period: 2104 (now the last row) keeps 30506,
# every other period (now starting with 2001 on the first row) is 35200.
$ws.Range("F16").Value = 35200
$ws.Range("F31").Value = 30506
